$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text in the source data.
# Force text number format so COM does not coerce them into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281.93"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.241"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06154"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.575"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.479"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8159"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1630"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08270"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03536"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03149"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09141"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.734"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04656"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006469"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001066"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.817"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.322"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3373"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04644"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007109"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1100"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003397"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01138"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006219"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002937"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001898"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01239"

# Row reshuffle + text updates for columns B, C, E (rows 41-43).
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
